$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 236 - this shifts existing rows 236:246 down to 237:247
# and automatically extends the sheet dimension to A1:R247.
$ws.Rows.Item(236).Insert()

# Populate the newly inserted row 236 with the new weekly price record
# (mirrors the layout of the surrounding "Poroto verde" rows).
$ws.Range("A236").Value = 8
$ws.Range("B236").Value = "Terminal La Palmera de La Serena"
$ws.Range("C236").Value = "Coquimbo"
$ws.Range("D236").Value = 44753
$ws.Range("E236").Value = 4
$ws.Range("F236").Value = 100112031
$ws.Range("G236").Value = "Poroto verde"
$ws.Range("H236").Value = "Magnum"
$ws.Range("I236").Value = "Primera"
$ws.Range("J236").Value = 440
$ws.Range("K236").Value = 35000
$ws.Range("L236").Value = 36000
$ws.Range("M236").Value = 35500
$ws.Range("N236").Value = "$/malla 25 kilos"
$ws.Range("O236").Value = "Perú"
$ws.Range("P236").Value = 1420
$ws.Range("Q236").Value = 25
$ws.Range("R236").Value = "Hortaliza"
